# Add javadoc and reporting
# Extend the "Locations" sheet with additional benchmark rows (rows 4-7),
# mirroring the data already present on the "Rewards" sheet, and correct
# a reported value on the "Rewards" sheet (C7: 1 -> 0).

$wb = $excel.ActiveWorkbook

# --- Locations sheet: append rows 4-7 ---
$ws = $wb.Worksheets.Item("Locations")

$newRows = @(
    @(5000.0,   9.0,   900.0),
    @(10000.0,  20.0,  900.0),
    @(50000.0,  105.0, 900.0),
    @(100000.0, 212.0, 900.0)
)

# Seed the new rows with the formatting of the last existing data row
# (style index 1) by copying it down, then overwrite the values.
$r = 4
foreach ($row in $newRows) {
    $ws.Range("B3:D3").Copy()
    $ws.Range("B" + $r + ":D" + $r).PasteSpecial()

    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]

    $r++
}

# --- Rewards sheet: fix reported value in C7 ---
$wsRewards = $wb.Worksheets.Item("Rewards")
$wsRewards.Cells.Item(7, 3).Value = 0.0
